{"js": "const NEW_VALUES = [\"1+4=5\", \"3+38=41\", \"17-9=8\", \"33+4=37\", \"30+47=77\", \"95-11=84\", \"18+9=27\", \"30-12=18\", \"27+25=52\", \"5+92=97\", \"65-62=3\", \"75-3=72\", \"6+93=99\", \"67-14=53\", \"35+38=73\", \"46-34=12\", \"80-3=77\", \"46-35=11\", \"66-18=48\", \"5+61=66\", \"60+3=63\", \"22-15=7\", \"20-8=12\", \"58-38=20\", \"18+25=43\", \"6+83=89\", \"7+58=65\", \"11+42=53\", \"86-69=17\", \"94-34=60\", \"83-41=42\", \"82-34=48\", \"58+21=79\", \"55+27=82\", \"9+50=59\", \"22+44=66\", \"90-34=56\", \"79-44=35\", \"86-11=75\", \"11+29=40\", \"52-49=3\", \"6-4=2\", \"17+75=92\", \"53-38=15\", \"92+6=98\", \"20+34=54\", \"84-2=82\", \"85-49=36\", \"89-84=5\", \"28+69=97\", \"21+67=88\", \"59-50=9\", \"17+55=72\", \"99-86=13\", \"18-5=13\", \"90-2=88\", \"98-15=83\", \"20+24=44\", \"0+71=71\", \"2+72=74\", \"69+25=94\", \"21+54=75\", \"67-31=36\", \"31-5=26\", \"93-14=79\", \"71-42=29\", \"99-87=12\", \"39+22=61\", \"15+79=94\", \"74-14=60\", \"31+38=69\", \"33+41=74\", \"91-18=73\", \"88-36=52\", \"42+21=63\", \"29+60=89\", \"55-37=18\", \"54-0=54\", \"44+20=64\", \"80-57=23\", \"6+1=7\", \"3+65=68\", \"84-61=23\", \"44+19=63\", \"44-31=13\", \"92-63=29\", \"72-35=37\", \"83-80=3\", \"64+15=79\", \"78+11=89\", \"20+13=33\", \"89-32=57\", \"79-66=13\", \"79-11=68\", \"43+14=57\", \"6+76=82\", \"6+73=79\", \"16-2=14\", \"87+2=89\", \"50-49=1\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body\");\n}\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst ROWS = table.rowCount;\nconst COLS = 5;\n\nif (ROWS * COLS !== NEW_VALUES.length) {\n  throw new Error(\"Unexpected table size: \" + ROWS + \"x\" + COLS + \" != \" + NEW_VALUES.length);\n}\n\n// Replace the text run-by-run, cell-by-cell, in row-major order (top-to-bottom,\n// left-to-right) \u2014 the same order the 100 <w:t> values appear in the document.\n// Using the paragraph's Range (rather than the cell body) preserves the\n// existing run/paragraph formatting (font, size, justification) instead of\n// collapsing it to defaults.\nlet k = 0;\nfor (let r = 0; r < ROWS; r++) {\n  for (let c = 0; c < COLS; c++) {\n    const cell = table.getCell(r, c);\n    const para = cell.body.paragraphs.getFirst();\n    para.getRange().insertText(NEW_VALUES[k], Word.InsertLocation.replace);\n    k++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The 100 replacement \"addition/subtraction within 100\" answers, in the same\n# row-major order (top-to-bottom, left-to-right) as the 20x5 table's cells.\n$NewValues = @(\"1+4=5\",\"3+38=41\",\"17-9=8\",\"33+4=37\",\"30+47=77\",\"95-11=84\",\"18+9=27\",\"30-12=18\",\"27+25=52\",\"5+92=97\",\"65-62=3\",\"75-3=72\",\"6+93=99\",\"67-14=53\",\"35+38=73\",\"46-34=12\",\"80-3=77\",\"46-35=11\",\"66-18=48\",\"5+61=66\",\"60+3=63\",\"22-15=7\",\"20-8=12\",\"58-38=20\",\"18+25=43\",\"6+83=89\",\"7+58=65\",\"11+42=53\",\"86-69=17\",\"94-34=60\",\"83-41=42\",\"82-34=48\",\"58+21=79\",\"55+27=82\",\"9+50=59\",\"22+44=66\",\"90-34=56\",\"79-44=35\",\"86-11=75\",\"11+29=40\",\"52-49=3\",\"6-4=2\",\"17+75=92\",\"53-38=15\",\"92+6=98\",\"20+34=54\",\"84-2=82\",\"85-49=36\",\"89-84=5\",\"28+69=97\",\"21+67=88\",\"59-50=9\",\"17+55=72\",\"99-86=13\",\"18-5=13\",\"90-2=88\",\"98-15=83\",\"20+24=44\",\"0+71=71\",\"2+72=74\",\"69+25=94\",\"21+54=75\",\"67-31=36\",\"31-5=26\",\"93-14=79\",\"71-42=29\",\"99-87=12\",\"39+22=61\",\"15+79=94\",\"74-14=60\",\"31+38=69\",\"33+41=74\",\"91-18=73\",\"88-36=52\",\"42+21=63\",\"29+60=89\",\"55-37=18\",\"54-0=54\",\"44+20=64\",\"80-57=23\",\"6+1=7\",\"3+65=68\",\"84-61=23\",\"44+19=63\",\"44-31=13\",\"92-63=29\",\"72-35=37\",\"83-80=3\",\"64+15=79\",\"78+11=89\",\"20+13=33\",\"89-32=57\",\"79-66=13\",\"79-11=68\",\"43+14=57\",\"6+76=82\",\"6+73=79\",\"16-2=14\",\"87+2=89\",\"50-49=1\")\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$rows = $table.Rows.Count\n$cols = 5\n\nif (($rows * $cols) -ne $NewValues.Count) {\n    throw \"Unexpected table size: $rows x $cols != $($NewValues.Count)\"\n}\n\n# Setting Cell.Range.Text (rather than replacing the whole cell) keeps the\n# existing run/paragraph formatting (font, size, justification) intact and\n# only swaps the visible text, matching the diff exactly.\n$k = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $table.Cell($r, $c).Range.Text = $NewValues[$k]\n        $k++\n    }\n}\n"}
